# Regional Availability Factor workbook update ("updated 4.0 files and mdl")
#
# 1. About sheet: bump the "last updated" date stamp (C1) from 2024-03-15 to 2024-03-28
# 2. RAF-capacity sheet: raise the capacity-credit multiplier for the two
#    hydrogen rows (hydrogen combustion turbine / hydrogen combined cycle)
#    from 0.3 to 1
# 3. Move the user's active view from RAF-generation to RAF-capacity, scrolled
#    down toward the bottom of the sheet with a reduced zoom level, selecting B25

$wb = $excel.ActiveWorkbook

# --- About sheet: update the date in C1 -------------------------------------------------
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Range("C1").Value = 45379

# --- RAF-capacity sheet: update the hydrogen capacity-credit multipliers ----------------
$wsCapacity = $wb.Worksheets.Item("RAF-capacity")
$wsCapacity.Range("B24").Value = 1
$wsCapacity.Range("B25").Value = 1

# Narrow column A slightly to better fit the shorter labels on this sheet
$wsCapacity.Columns.Item(1).ColumnWidth = 28.125

# --- Switch the active sheet/view to RAF-capacity, scrolled and zoomed --------------------
$wsCapacity.Activate()
$wsCapacity.Range("B25").Select()
$excel.ActiveWindow.Zoom = 80
$excel.ActiveWindow.ScrollRow = 14
$excel.ActiveWindow.ScrollColumn = 1
